# Updated BGR model - 2025-08-17 18:52
# Re-rank the "lcoe_class" (column P) values for a handful of resource
# classes on the "solar" and "wind" sheets.

$wb = $excel.ActiveWorkbook

# --- solar sheet: CF class spv-BGR_16 ------------------------------------
# old ranks: c3=3, c4=4, c2=2  ->  new ranks: c3=4 (row4), c4=2 (row5), c2=3 (row6)
$wsSolar = $wb.Worksheets.Item("solar")
$wsSolar.Range("P4").Value = 4
$wsSolar.Range("P5").Value = 2
$wsSolar.Range("P6").Value = 3

# --- wind sheet: CF classes won-BGR_25, won-BGR_24, won-BGR_17 ----------
$wsWind = $wb.Worksheets.Item("wind")

# won-BGR_25: c4 (row13) <-> c5 (row14)
$wsWind.Range("P13").Value = 5
$wsWind.Range("P14").Value = 4

# won-BGR_25: c1 (row16) <-> c2 (row17)
$wsWind.Range("P16").Value = 2
$wsWind.Range("P17").Value = 1

# won-BGR_24: c1 (row19) <-> c3 (row20)
$wsWind.Range("P19").Value = 3
$wsWind.Range("P20").Value = 1

# won-BGR_17: c1 (row47) <-> c2 (row48)
$wsWind.Range("P47").Value = 2
$wsWind.Range("P48").Value = 1
